$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.348.08"
$ws.Range("E2").Value = "  +2.46%  "
$ws.Range("D3").Value = "2.475.18"
$ws.Range("E3").Value = "  +2.56%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.29"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.04"
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("D9").Value = "2.475.80"
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("E10").Value = "  +2.49%  "
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.13"
$ws.Range("E14").Value = "  +9.46%  "
$ws.Range("E15").Value = "  +3.22%  "
$ws.Range("D16").Value = "2.918.60"
$ws.Range("E16").Value = "  +2.93%  "
$ws.Range("D17").Value = "63.263.05"
$ws.Range("E17").Value = "  +2.68%  "
$ws.Range("D18").Value = "2.465.44"
$ws.Range("E18").Value = "  +2.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.93"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.12"
$ws.Range("E20").Value = "  +3.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.36"
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("E22").Value = "  +11.07%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("B25").Value = "Bittensor"
$ws.Range("C25").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "678.84"
$ws.Range("E25").Value = "  +9.52%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "66.57"
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.23"
$ws.Range("E27").Value = "  +9.98%  "
$ws.Range("E28").Value = "  +5.08%  "
$ws.Range("D29").Value = "2.597.55"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.01"
$ws.Range("E30").Value = "  +1.54%  "
$ws.Range("E31").Value = "  +3.34%  "
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("E33").Value = "  +4.03%  "
$ws.Range("E34").Value = "  +3.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  +6.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("E37").Value = "  +3.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.57"
$ws.Range("E38").Value = "  +3.93%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "153.78"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.374"
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.89"
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.78"
$ws.Range("E42").Value = "  +6.42%  "
$ws.Range("E43").Value = "  +3.87%  "
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "0.0₆0300"
$ws.Range("E46").Value = "  +5.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.16"
$ws.Range("E47").Value = "  +27.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "146.93"
$ws.Range("E48").Value = "  +3.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.64"
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.87"
$ws.Range("E50").Value = "  +4.15%  "
$ws.Range("E51").Value = "  +2.09%  "
